$d = $word.ActiveDocument

$find = "Campaign Dates that use Perseus constellation 2022: January 16-25, November 7-16, December 6-15"
$replace = " 2022 Campaign Dates that use Perseus constellation: January 16-25, November 7-16, December 6-15"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
